# Rename the existing sheet to "FirstSheet" and add a second sheet
# ("SecondSheet") that is a full copy of it (same layout/styles/columns),
# except the host-name value in the second row is different.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)

# Duplicate FirstSheet right after itself - this clones all data, column
# widths, row heights and per-cell styles exactly.
$ws1.Copy([System.Reflection.Missing]::Value, $ws1)

# Now rename the two sheets.
$ws1.Name = "FirstSheet"
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "SecondSheet"

# The only content difference on the new sheet: A2 gets its own value.
$ws2.Range("A2").Value = "mb01第2页啦"

# Keep FirstSheet as the active/selected sheet, same as the original file.
$ws1.Activate()
